$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing A14 text from "how" to "Testing git ignore"
$ws.Range("A14").Value = "Testing git ignore"

# Add new row A15 with "Testing git ignore2"
$ws.Range("A15").Value = "Testing git ignore2"

# Move the active selection to A16, matching the post-edit sheet view
$ws.Activate()
$ws.Range("A16").Select()
